$wb = $excel.ActiveWorkbook

# The "approver" name on the approvalDetails sheet was updated.
$ws = $wb.Worksheets.Item("approvalDetails")
$ws.Range("D2").Value = "K.Ramakrishna-ADM_Senior Assistant_7"

# approvalDetails becomes the active/selected sheet (was challanHeaderDetails),
# with the selection moved to D2 (was E6).
$ws.Activate() | Out-Null
$ws.Range("D2").Select() | Out-Null
